$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly data row was added to the top of the price history block
# (row 322). Insert a blank row there, which pushes the existing rows
# 322:359 down to 323:360 (carrying their values/formatting with them),
# then populate the newly inserted row with this week's figures.
$ws.Rows(322).Insert()

$ws.Cells.Item(322, 1).Value = 9
$ws.Cells.Item(322, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(322, 3).Value = "Metropolitana"
$ws.Cells.Item(322, 4).Value = 45124
$ws.Cells.Item(322, 5).Value = 13
$ws.Cells.Item(322, 6).Value = 100112026
$ws.Cells.Item(322, 7).Value = "Haba"
$ws.Cells.Item(322, 8).Value = "Sin especificar"
$ws.Cells.Item(322, 9).Value = "Primera"
$ws.Cells.Item(322, 10).Value = 52
$ws.Cells.Item(322, 11).Value = 14000
$ws.Cells.Item(322, 12).Value = 16000
$ws.Cells.Item(322, 13).Value = 15000
$ws.Cells.Item(322, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(322, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(322, 16).Value = 600
$ws.Cells.Item(322, 17).Value = 25
$ws.Cells.Item(322, 18).Value = "Hortaliza"
